$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "محمود أحمد شوقي إبراهيم"

$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:ma7moud.a.mojahed@gmail.com", "", "", "ma7moud.a.mojahed@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/M-Mojahed/OP_Source", "", "", "https://github.com/M-Mojahed/OP_Source")

[void]$ws.Range("C3").Select()
